$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-6 from 45243 to 45244
$ws.Range("C2:C6").Value = 45244
